$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EVCRSbRIC")

# The original column K held "ISIC 20T21" covering both ISIC 20 and ISIC 21.
# Split it into two columns: K = "ISIC 20", and a newly inserted L = "ISIC 21".
# Insert a new column before the old column L (i.e. after K), shifting
# everything from L onward one column to the right.
$ws.Range("L1").EntireColumn.Insert()

# Rename the header that used to say "ISIC 20T21" to just "ISIC 20".
$ws.Range("K1").Value = "ISIC 20"

# Populate the header of the newly inserted column with "ISIC 21".
$ws.Range("L1").Value = "ISIC 21"

# The data row below (row 2) keeps a 0 value in the newly-inserted column,
# matching the pattern used by the other ISIC columns.
$ws.Range("L2").Value = 0
